# daily_increment_positif.xlsx — "update sampai 26 april"
# New daily data through 2020-04-26 (serial 43947): row 37 gets two
# corrected values (highlighted), and four brand-new rows (38-41) are
# appended for 2020-04-23 .. 2020-04-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$HIGHLIGHT = 16776960   # RGB(0,255,255) cyan — matches the workbook's existing "corrected value" fill

# --- Row 37 (2020-04-22) corrections ---
$ws.Cells.Item(37, 8).Value = 120          # H37: 119 -> 120
$ws.Cells.Item(37, 8).Interior.Color = $HIGHLIGHT
$ws.Cells.Item(37, 36).Value = -1          # AJ37: 0 -> -1
$ws.Cells.Item(37, 36).Interior.Color = $HIGHLIGHT

# --- Row 38 (2020-04-23) ---
$ws.Rows.Item(38).RowHeight = 13.5
$ws.Cells.Item(38, 1).Value = 43944
$ws.Cells.Item(38, 1).NumberFormat = "yyyy-mmm-dd"
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(38, 3).Value = 15
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 1
$ws.Cells.Item(38, 8).Value = 130
$ws.Cells.Item(38, 8).Interior.Color = $HIGHLIGHT
$ws.Cells.Item(38, 9).Value = 1
$ws.Cells.Item(38, 10).Value = 25
$ws.Cells.Item(38, 10).Interior.Color = $HIGHLIGHT
$ws.Cells.Item(38, 11).Value = 59
$ws.Cells.Item(38, 12).Value = 26
$ws.Cells.Item(38, 13).Value = 19
$ws.Cells.Item(38, 14).Value = 5
$ws.Cells.Item(38, 15).Value = 1
$ws.Cells.Item(38, 16).Value = 7
$ws.Cells.Item(38, 17).Value = 0
$ws.Cells.Item(38, 18).Value = 2
$ws.Cells.Item(38, 19).Value = 7
$ws.Cells.Item(38, 20).Value = 4
$ws.Cells.Item(38, 21).Value = 5
$ws.Cells.Item(38, 22).Value = 11
$ws.Cells.Item(38, 23).Value = 2
$ws.Cells.Item(38, 24).Value = 0
$ws.Cells.Item(38, 25).Value = 10
$ws.Cells.Item(38, 26).Value = 0
$ws.Cells.Item(38, 27).Value = 11
$ws.Cells.Item(38, 28).Value = 1
$ws.Cells.Item(38, 29).Value = 2
$ws.Cells.Item(38, 30).Value = 0
$ws.Cells.Item(38, 31).Value = 5
$ws.Cells.Item(38, 32).Value = 7
$ws.Cells.Item(38, 33).Value = 0
$ws.Cells.Item(38, 34).Value = 0
$ws.Cells.Item(38, 35).Value = 0
$ws.Cells.Item(38, 36).Value = 0

# --- Row 39 (2020-04-24) ---
$ws.Rows.Item(39).RowHeight = 13.5
$ws.Cells.Item(39, 1).Value = 43945
$ws.Cells.Item(39, 1).NumberFormat = "yyyy-mmm-dd"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = 10
$ws.Cells.Item(39, 4).Value = 22
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = 85
$ws.Cells.Item(39, 9).Value = 4
$ws.Cells.Item(39, 10).Value = 75
$ws.Cells.Item(39, 11).Value = 37
$ws.Cells.Item(39, 12).Value = 26
$ws.Cells.Item(39, 13).Value = 0
$ws.Cells.Item(39, 14).Value = 11
$ws.Cells.Item(39, 15).Value = 11
$ws.Cells.Item(39, 16).Value = 18
$ws.Cells.Item(39, 17).Value = 0
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 19).Value = 38
$ws.Cells.Item(39, 20).Value = 13
$ws.Cells.Item(39, 21).Value = 10
$ws.Cells.Item(39, 22).Value = 5
$ws.Cells.Item(39, 23).Value = 1
$ws.Cells.Item(39, 24).Value = 4
$ws.Cells.Item(39, 25).Value = 23
$ws.Cells.Item(39, 26).Value = 3
$ws.Cells.Item(39, 27).Value = 0
$ws.Cells.Item(39, 28).Value = 0
$ws.Cells.Item(39, 29).Value = 0
$ws.Cells.Item(39, 30).Value = 0
$ws.Cells.Item(39, 31).Value = 2
$ws.Cells.Item(39, 32).Value = 6
$ws.Cells.Item(39, 33).Value = 25
$ws.Cells.Item(39, 34).Value = 0
$ws.Cells.Item(39, 35).Value = 5
$ws.Cells.Item(39, 36).Value = 0

# --- Row 40 (2020-04-25) ---
$ws.Rows.Item(40).RowHeight = 13.5
$ws.Cells.Item(40, 1).Value = 43946
$ws.Cells.Item(40, 1).NumberFormat = "yyyy-mmm-dd"
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = 6
$ws.Cells.Item(40, 4).Value = 11
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 22
$ws.Cells.Item(40, 8).Value = 85
$ws.Cells.Item(40, 9).Value = 3
$ws.Cells.Item(40, 10).Value = 45
$ws.Cells.Item(40, 11).Value = 46
$ws.Cells.Item(40, 12).Value = 80
$ws.Cells.Item(40, 13).Value = 1
$ws.Cells.Item(40, 14).Value = 12
$ws.Cells.Item(40, 15).Value = 6
$ws.Cells.Item(40, 16).Value = 14
$ws.Cells.Item(40, 17).Value = 6
$ws.Cells.Item(40, 18).Value = 0
$ws.Cells.Item(40, 19).Value = 27
$ws.Cells.Item(40, 20).Value = 13
$ws.Cells.Item(40, 21).Value = 1
$ws.Cells.Item(40, 22).Value = 0
$ws.Cells.Item(40, 23).Value = 9
$ws.Cells.Item(40, 24).Value = 0
$ws.Cells.Item(40, 25).Value = 12
$ws.Cells.Item(40, 26).Value = 4
$ws.Cells.Item(40, 27).Value = 0
$ws.Cells.Item(40, 28).Value = 22
$ws.Cells.Item(40, 29).Value = 0
$ws.Cells.Item(40, 30).Value = 5
$ws.Cells.Item(40, 31).Value = 1
$ws.Cells.Item(40, 32).Value = 0
$ws.Cells.Item(40, 33).Value = 2
$ws.Cells.Item(40, 34).Value = 0
$ws.Cells.Item(40, 35).Value = 2
$ws.Cells.Item(40, 36).Value = 0

# --- Row 41 (2020-04-26) ---
$ws.Rows.Item(41).RowHeight = 13.5
$ws.Cells.Item(41, 1).Value = 43947
$ws.Cells.Item(41, 1).NumberFormat = "yyyy-mmm-dd"
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 3
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 1
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 3
$ws.Cells.Item(41, 8).Value = 114
$ws.Cells.Item(41, 9).Value = 11
$ws.Cells.Item(41, 10).Value = 5
$ws.Cells.Item(41, 11).Value = 28
$ws.Cells.Item(41, 12).Value = 15
$ws.Cells.Item(41, 13).Value = 0
$ws.Cells.Item(41, 14).Value = 8
$ws.Cells.Item(41, 15).Value = 4
$ws.Cells.Item(41, 16).Value = 0
$ws.Cells.Item(41, 17).Value = 6
$ws.Cells.Item(41, 18).Value = 2
$ws.Cells.Item(41, 19).Value = 15
$ws.Cells.Item(41, 20).Value = 10
$ws.Cells.Item(41, 21).Value = 5
$ws.Cells.Item(41, 22).Value = 4
$ws.Cells.Item(41, 23).Value = 6
$ws.Cells.Item(41, 24).Value = 4
$ws.Cells.Item(41, 25).Value = 8
$ws.Cells.Item(41, 26).Value = 0
$ws.Cells.Item(41, 27).Value = 4
$ws.Cells.Item(41, 28).Value = 1
$ws.Cells.Item(41, 29).Value = 12
$ws.Cells.Item(41, 30).Value = 0
$ws.Cells.Item(41, 31).Value = 0
$ws.Cells.Item(41, 32).Value = 5
$ws.Cells.Item(41, 33).Value = 0
$ws.Cells.Item(41, 34).Value = 0
$ws.Cells.Item(41, 35).Value = 1
$ws.Cells.Item(41, 36).Value = 0

# --- Cell comments (author notes for the highlighted corrections / new entries) ---
$ws.Cells.Item(37, 8).AddComment("119")
$ws.Cells.Item(37, 36).AddComment("0")
$ws.Cells.Item(38, 8).AddComment("133")
$ws.Cells.Item(38, 10).AddComment("22")

# --- View: scroll down to the new rows and select the bottom-right corner cell ---
$ws.Activate()
$ws.Range("AK41").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
